# Generate Report for Handoff
#
# The "b.md" file has finished handoff processing: its status flips from
# "Handed back: in sync with en-US" to "Ready for handoff", a new handoff
# package (b.<hash>.<locale>.xlf) was generated, content duplication is no
# longer flagged, and (for each locale) an error is now reported because the
# handback file on record is stale compared to the newly generated b.md
# handoff.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ddbad2cfa3d31c6ea564c22a29ac5df628619204/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ed14a41323f896b0922b3f4aa58fb09d5ea7e213/e2e/b.md."

# ---------------------------------------------------------------------
# Overview sheet: row 3 is the b.md roll-up row.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-30 20:42:34"

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 is the b.md detail row.
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("F3").Style = "Normal"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-30 20:42:29"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Range("P1").EntireColumn.ColumnWidth = 39.1

# ---------------------------------------------------------------------
# de-de sheet: row 3 is the b.md detail row.
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("F3").Style = "Normal"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-30 20:42:34"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Range("P1").EntireColumn.ColumnWidth = 39.1
